$wb = $excel.ActiveWorkbook

$colIndex = @{ "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11 }

function Set-SheetData {
    param($ws, $data)
    foreach ($r in $data.Keys) {
        $row = $data[$r]
        foreach ($c in $row.Keys) {
            $ws.Cells.Item([int]$r, $colIndex[$c]).Value = [double]$row[$c]
        }
    }
}

# ---- Update annualised_return sheet ----
$data_sheet12 = @{
    2 = @{ "B"=1; "C"=1; "D"=0.0003; "E"=0.0002; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    3 = @{ "B"=1; "C"=1; "D"=0; "E"=0; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    4 = @{ "B"=0.0003; "C"=0; "D"=1; "E"=1; "F"=0.0002; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    5 = @{ "B"=0.0002; "C"=0; "D"=1; "E"=1; "F"=0.0002; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    6 = @{ "B"=1; "C"=1; "D"=0.0002; "E"=0.0002; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    7 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=1; "I"=1; "J"=1; "K"=1 }
    8 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=1; "I"=1; "J"=1; "K"=0.307 }
    9 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=1; "I"=1; "J"=1; "K"=1 }
    10 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=1; "I"=1; "J"=1; "K"=1 }
    11 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=0.307; "I"=1; "J"=1; "K"=1 }
}
$wsAnn = $wb.Worksheets.Item("annualised_return")
Set-SheetData $wsAnn $data_sheet12

# ---- Update mean_period_return sheet (identical data) ----
$wsMean = $wb.Worksheets.Item("mean_period_return")
Set-SheetData $wsMean $data_sheet12

# ---- Update sharpe_annualized sheet ----
$data_sheet3 = @{
    2 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0.1312; "I"=0.429; "J"=1; "K"=1 }
    3 = @{ "B"=0; "C"=1; "D"=0.0038; "E"=0.0036; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    4 = @{ "B"=0; "C"=0.0038; "D"=1; "E"=1; "F"=0.0001; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    5 = @{ "B"=0; "C"=0.0036; "D"=1; "E"=1; "F"=0.0001; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    6 = @{ "B"=0; "C"=1; "D"=0.0001; "E"=0.0001; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    7 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=0; "I"=0; "J"=0; "K"=0 }
    8 = @{ "B"=0.1312; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0 }
    9 = @{ "B"=0.429; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0003 }
    10 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0221 }
    11 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0; "I"=0.0003; "J"=0.0221; "K"=1 }
}
$wsSharpeAnn = $wb.Worksheets.Item("sharpe_annualized")
Set-SheetData $wsSharpeAnn $data_sheet3

# ---- Insert new sharpe_period sheet before VaR (copy formatting from sharpe_annualized) ----
$wsVaR = $wb.Worksheets.Item("VaR")
$wsSharpeAnn.Copy($wsVaR)
# Re-fetch references by name: after Copy, old $wsVaR variable now points to the newly inserted sheet
$wsSharpePeriod = $wb.Worksheets.Item("sharpe_annualized (2)")
$wsSharpePeriod.Name = "sharpe_period"
$data_sharpe_period = @{
    2 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    3 = @{ "B"=0; "C"=1; "D"=0.7118; "E"=0.6922; "F"=0.1089; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    4 = @{ "B"=0; "C"=0.7118; "D"=1; "E"=1; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    5 = @{ "B"=0; "C"=0.6922; "D"=1; "E"=1; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    6 = @{ "B"=0; "C"=0.1089; "D"=1; "E"=1; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    7 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=0; "I"=0; "J"=0; "K"=0 }
    8 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.1005 }
    9 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0256 }
    10 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0616 }
    11 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0.1005; "I"=0.0256; "J"=0.0616; "K"=1 }
}
Set-SheetData $wsSharpePeriod $data_sharpe_period

# ---- Update VaR sheet (re-fetch reference by name, since it shifted index after insert) ----
$wsVaR = $wb.Worksheets.Item("VaR")
$data_sheet4 = @{
    2 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0.0175; "I"=0.0672; "J"=1; "K"=1 }
    3 = @{ "B"=0; "C"=1; "D"=0.0048; "E"=0.0046; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    4 = @{ "B"=0; "C"=0.0048; "D"=1; "E"=1; "F"=0.0001; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    5 = @{ "B"=0; "C"=0.0046; "D"=1; "E"=1; "F"=0.0001; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    6 = @{ "B"=0; "C"=1; "D"=0.0001; "E"=0.0001; "F"=1; "G"=0; "H"=0; "I"=0; "J"=0; "K"=0 }
    7 = @{ "B"=0; "C"=0; "D"=0; "E"=0; "F"=0; "G"=1; "H"=0; "I"=0; "J"=0; "K"=0 }
    8 = @{ "B"=0.0175; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0 }
    9 = @{ "B"=0.0672; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0002 }
    10 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=1; "I"=1; "J"=1; "K"=0.0191 }
    11 = @{ "B"=1; "C"=0; "D"=0; "E"=0; "F"=0; "G"=0; "H"=0; "I"=0.0002; "J"=0.0191; "K"=1 }
}
Set-SheetData $wsVaR $data_sheet4

$wsAnn = $wb.Worksheets.Item("annualised_return")
$wsAnn.Activate()

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Index $s.Name
}